# Adds a new worksheet "2025-09-10" after "2025-09-09" with the
# day's ranking data (header + 50 rows), matching the style of the
# other date sheets (bold/bordered/centered header row).

$wb = $excel.ActiveWorkbook

$anchor = $wb.Worksheets.Item('2025-09-09')
$ws = $wb.Worksheets.Add($null, $anchor)
$ws.Name = '2025-09-10'

# Reuse the exact header-row cell style from the previous date sheet
# (bold font, thin border, centered/top-aligned) instead of building a
# fresh style, so the new sheet shares the same style index as the rest.
$anchor.Range('A1:D1').Copy()
$ws.Range('A1:D1').PasteSpecial(-4122)
$excel.CutCopyMode = $false

$data = New-Object 'object[,]' 51,4
$data[0,0] = 'rank'
$data[0,1] = 'title'
$data[0,2] = 'author'
$data[0,3] = 'latest_episode'

$data[1,0] = 1
$data[1,1] = '転生コロシアム～最弱スキルで最強の女たちを攻略して奴隷ハーレム作ります～'
$data[1,2] = 'zunta(作画) はらわたさいぞう(原作)'
$data[1,3] = '第32話：思考を奪う②'
$data[2,0] = 2
$data[2,1] = '時間停止勇者―余命３日の設定じゃ世界を救うには短すぎる―'
$data[2,2] = '光永康則'
$data[2,3] = '第６９話『岩鬼停止』②'
$data[3,0] = 3
$data[3,1] = '地元のいじめっ子達に仕返ししようとしたら、別の戦いが始まった。'
$data[3,2] = 'マツモトケンゴ'
$data[3,3] = '第64話 更衣室の戦いが始まった（２）'
$data[4,0] = 4
$data[4,1] = '望まぬ不死の冒険者'
$data[4,2] = '中曽根ハイジ（漫画） 丘野 優（原作） じゃいあん（キャラクター原案）'
$data[4,3] = '第60話　異国の観光客'
$data[5,0] = 5
$data[5,1] = '世界最強の魔女、始めました 〜私だけ『攻略サイト』を見れる世界で自由に生きます〜'
$data[5,2] = '戸賀 環 坂木持丸 riritto'
$data[5,3] = '第52話②　最強の武器を手に入れてみた'
$data[6,0] = 6
$data[6,1] = '配信に致命的に向いていない女の子が迷宮で黙々と人助けする配信'
$data[6,2] = '下田将也(漫画) 佐藤悪糖(原作) 福きつね(キャラクター原案)'
$data[6,3] = '第2話後編'
$data[7,0] = 7
$data[7,1] = '乙女ゲー世界はモブに厳しい世界です【共和国編】'
$data[7,2] = '三嶋与夢(原作) 行々狸(作画) 孟達(キャラクター原案) マツリセイシロウ(構成) FTops(制作)'
$data[7,3] = '第2話-1'
$data[8,0] = 8
$data[8,1] = 'ラスボス討伐後に始める二周目冒険者ライフ はじまりの街でワケあり美少女たちがめちゃくちゃ懐いてきます'
$data[8,2] = '鬼麻正明(漫画) 朱月十話(原作) ファルまろ(キャラ原案)'
$data[8,3] = '第5話-1'
$data[9,0] = 9
$data[9,1] = '絶対死なないステラ姫'
$data[9,2] = '光永康則 大高稲'
$data[9,3] = '第１５話　絶対指名手配されない（３）'
$data[10,0] = 10
$data[10,1] = '女子高生の無駄づかい'
$data[10,2] = 'ビーノ(著者)'
$data[10,3] = '第135話　きずな'
$data[11,0] = 11
$data[11,1] = 'ひとりぼっちの異世界攻略'
$data[11,2] = 'びび（漫画） 五示正司（原作）'
$data[11,3] = '第235話　色々あるって言ってるじゃん'
$data[12,0] = 12
$data[12,1] = '序盤で死ぬ最強のサブキャラに転生したので、ゲーム知識で無双する'
$data[12,2] = '作画：マエD 原作：新人'
$data[12,3] = '第6話(4)'
$data[13,0] = 13
$data[13,1] = '江戸前エルフ'
$data[13,2] = '樋口彰彦'
$data[13,3] = '#122'
$data[14,0] = 14
$data[14,1] = 'ホームセンターごと呼び出された私の大迷宮リノベーション！'
$data[14,2] = 'ばたっち(漫画) 星崎崑(原作) 志田(キャラクター原案)'
$data[14,3] = '第5話前編'
$data[15,0] = 15
$data[15,1] = 'この素晴らしい世界に祝福を！'
$data[15,2] = '渡真仁(作画) 三嶋くろね(キャラクター原案) 暁なつめ(原作)'
$data[15,3] = '第131話-2　この冒険者達と共に原点回帰を！②'
$data[16,0] = 16
$data[16,1] = 'ある日、惰眠を貪っていたら一族から追放されて森に捨てられました そのまま寝てたら周りが勝手に魔物の国を作ってたけど、私は気にせず今日も眠ります　コミック版'
$data[16,2] = '漫画/伊草さゆ 原作/白波ハクア キャラクター原案/まさよ'
$data[16,3] = 'chapter55【29話①】'
$data[17,0] = 17
$data[17,1] = '俺にトラウマを与えた女子達がチラチラ見てくるけど、残念ですが手遅れです'
$data[17,2] = 'いちたか（漫画） 御堂ユラギ（原作） 緜（キャラクター原案）'
$data[17,3] = '第20話　母親'
$data[18,0] = 18
$data[18,1] = '無慈悲な悪役貴族に転生した僕は掌握魔法を駆使して魔法世界の頂点に立つ 〜ヒロインなんていないと諦めていたら向こうから勝手に寄ってきました〜'
$data[18,2] = '坂井オイ(漫画) びゃくし(原作) ファルまろ(キャラクター原案)'
$data[18,3] = '第7話-1'
$data[19,0] = 19
$data[19,1] = 'ダウナー系お姉さんに毎日カスの嘘を流し込まれる話'
$data[19,2] = '生倉のゑる(著者) はるばーど屋(原作者)'
$data[19,3] = '12話 おまけ'
$data[20,0] = 20
$data[20,1] = '傭兵団の料理番'
$data[20,2] = '梅木泰祐(漫画) 川井昂(原作) 四季童子(キャラクター原案)'
$data[20,3] = '第10話-1'
$data[21,0] = 21
$data[21,1] = 'ゲーム世界で魔物に転生してしまった俺、前世で推しだったヒロインを拾ってしまう'
$data[21,2] = '三部べべ(漫画) ねうしとら(原作)'
$data[21,3] = '第1話'
$data[22,0] = 22
$data[22,1] = '米原くんはつよつよギャルから離れられない'
$data[22,2] = '川村拓(著者)'
$data[22,3] = '第17話'
$data[23,0] = 23
$data[23,1] = '黒の召喚士'
$data[23,2] = '天羽 銀（漫画） 迷井豆腐（原作） 黒銀（DIGS）（キャラクター原案）'
$data[23,3] = '第148話　聖槍イクリプスⅨ'
$data[24,0] = 24
$data[24,1] = 'みだりに憑かせてはなりません'
$data[24,2] = '栗田あぐり(著者)'
$data[24,3] = '第10話①'
$data[25,0] = 25
$data[25,1] = '死ぬ運命にある悪役令嬢の兄に転生したので、妹を育てて未来を変えたいと思います　～世界最強はオレだけど、世界最カワは妹に違いない～'
$data[25,2] = '石見翔子(漫画） 泉里侑希（原作） タムラヨウ（キャラクター原案）'
$data[25,3] = '第5話　兄妹の約束（前編）'
$data[26,0] = 26
$data[26,1] = '追放されたチート付与魔術師は 気ままなセカンドライフを謳歌する。'
$data[26,2] = '六志麻あさ 業務用餅 kisui'
$data[26,3] = '第７１話'
$data[27,0] = 27
$data[27,1] = '蜘蛛ですが、なにか？'
$data[27,2] = 'かかし朝浩(著者) 馬場翁(原作) 輝竜司(キャラクター原案)'
$data[27,3] = '第76話その1'
$data[28,0] = 28
$data[28,1] = '美人女上司滝沢さん'
$data[28,2] = 'やんBARU(著者)'
$data[28,3] = '第202.5話'
$data[29,0] = 29
$data[29,1] = '実は俺、最強でした？'
$data[29,2] = '原作：澄守 彩 漫画：高橋 愛'
$data[29,3] = '第124話　ゴルドとナタリアとハルト'
$data[30,0] = 30
$data[30,1] = 'ラーメン大好き小泉さん'
$data[30,2] = '鳴見なる'
$data[30,3] = '24杯目 行列'
$data[31,0] = 31
$data[31,1] = '勇者パーティを追放された【スキルサポーター】、仲間のスキルを解放して最強に成り上がる'
$data[31,2] = '作画：なかお 原作：前田氏'
$data[31,3] = '第8話(2)'
$data[32,0] = 32
$data[32,1] = '悪役一家の奥方、死に戻りして心を入れ替える。'
$data[32,2] = '鏡(漫画) 丘野優(原作) TEDDY(キャラクター原案)'
$data[32,3] = '第33話②'
$data[33,0] = 33
$data[33,1] = '最強の少年聖騎士、転生者を狩る'
$data[33,2] = '作画：御塩 原作：宇奈木ユラ'
$data[33,3] = '第8話(2)'
$data[34,0] = 34
$data[34,1] = '二周目チートの転生魔導士～最強が1000年後に転生したら、人生余裕すぎました～'
$data[34,2] = '石後千鳥 鬱沢色素 りいちゅ'
$data[34,3] = '第32話　肝試し（後編）'
$data[35,0] = 35
$data[35,1] = 'オークの酒杯に祝福を'
$data[35,2] = 'かなどめはじめ'
$data[35,3] = '第48話　死神'
$data[36,0] = 36
$data[36,1] = '現代ダンジョンライフの続きは異世界オープンワールドで！'
$data[36,2] = '田中清久（漫画） しば犬部隊（原作） ひろせ（原作イラスト）'
$data[36,3] = '第25話　長い1日の始まり'
$data[37,0] = 37
$data[37,1] = '勇者パーティを追い出された器用貧乏　～パーティ事情で付与術士をやっていた剣士、万能へと至る～'
$data[37,2] = '漫画：よねぞう 原作：都神樹 キャラクター原案：きさらぎゆり'
$data[37,3] = '第５２話　暴走を止める器用貧乏（２）'
$data[38,0] = 38
$data[38,1] = '第七魔王子ジルバギアスの魔王傾国記'
$data[38,2] = '野井ニトラ（漫画） 甘木智彬（原作） 輝竜 司（キャラクター原案）'
$data[38,3] = '第18話　誤解です母上'
$data[39,0] = 39
$data[39,1] = '魔導具師ダリヤはうつむかない ～Dahliya Wilts No More～'
$data[39,2] = '漫画：住川惠 原作：甘岸久弥(｢魔導具師ダリヤはうつむかない ～今日から自由な職人ライフ～｣MFブックス刊) キャラクター原案：景、駒田ハチ'
$data[39,3] = '第47話 魔導具師とつながれたもの④'
$data[40,0] = 40
$data[40,1] = '迷宮狂走曲　～エロゲ世界なのにエロそっちのけでひたすら最強を目指すモブ転生者～'
$data[40,2] = 'ぱらボら（漫画） 宮迫宗一郎（原作） 灯（キャラクター原案）'
$data[40,3] = '第2話　「ＨＰ消費技ブッパ」は狂人の発想（前編）'
$data[41,0] = 41
$data[41,1] = '【パクパクですわ】追放されたお嬢様の『モンスターを食べるほど強くなる』スキルは、１食で１レベルアップする前代未聞の最強スキルでした。３日で人類最強になりましたわ～！'
$data[41,2] = '島知宏 音速炒飯 有都あらゆる'
$data[41,3] = '第２３食　巨大ヘビモンスターさん、パクパクですわ！（４）'
$data[42,0] = 42
$data[42,1] = '最凶貴族は死亡フラグを覆す'
$data[42,2] = '作画：sudekuma 原作：塚上'
$data[42,3] = '第8話(2)'
$data[43,0] = 43
$data[43,1] = '二度追放された冒険者、激レアスキル駆使して美少女軍団を育成中！　コミック版'
$data[43,2] = '漫画/青木千尋 原作/南野雪花'
$data[43,3] = 'chapter70【36話①】'
$data[44,0] = 44
$data[44,1] = '転生したら没落貴族だったので、【呪言】を極めて家族を救います'
$data[44,2] = '作画：アマセケイ 原作：メソポ・たみあ'
$data[44,3] = '第8話(2)'
$data[45,0] = 45
$data[45,1] = 'この冒険者、人類史最強です～外れスキル『鑑定』が『継承』に覚醒したので、数多の英雄たちの力を受け継ぎ無双する～'
$data[45,2] = '日之影ソラ みやけりく エシュアル'
$data[45,3] = '第29話②ダークエルフ救出作戦'
$data[46,0] = 46
$data[46,1] = '弱小国家の英雄王子　～最強の魔術師だけど、さっさと国出て自由に生きてぇぇ！～'
$data[46,2] = '友山アキ（漫画） 楓原こうた（原作） トモゼロ（キャラクター原案）'
$data[46,3] = '第4話　セリアの追憶（後編）'
$data[47,0] = 47
$data[47,1] = '十年目、帰還を諦めた転移者はいまさら主人公になる'
$data[47,2] = '原作：氷純（「十年目、帰還を諦めた転移者はいまさら主人公になる」MFブックス刊） 漫画：しゅーかま キャラクター原案：あんべよしろう'
$data[47,3] = '第１９話①'
$data[48,0] = 48
$data[48,1] = '無気力ニートな元神童、冒険者になる'
$data[48,2] = '緑茶こいめ（漫画） ぺもぺもさん（原作） 福きつね（原作イラスト）'
$data[48,3] = '第53話　元神童VS八星ヴァニール＆腹心ライオスⅡ'
$data[49,0] = 49
$data[49,1] = 'ブラックな騎士団の奴隷がホワイトな冒険者ギルドに引き抜かれてSランクになりました'
$data[49,2] = 'ハム梟（漫画） 寺王（原作） 由夜（キャラクター原案）'
$data[49,3] = '第44話　紅い女帝ⅩⅢ'
$data[50,0] = 50
$data[50,1] = '異世界でも無難に生きたい症候群'
$data[50,2] = '原作：安泰（一二三書房刊） 漫画：笹峰コウ キャラクター原案：ひたきゆう'
$data[50,3] = '第31話③'

$ws.Range('A1:D51').Value = $data

Write-Host "Worksheets: $($wb.Worksheets.Count); new sheet: $($ws.Name)"
